$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells E1, F1 (bold/centered/bordered header style like A1:D1)
$ws.Range("E1").Value = "avg_prompt_processing_rate_toks_per_sec"
$ws.Range("F1").Value = "avg_token_generation_rate_toks_per_sec"
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate E2:F58 with per-subject perf metrics
$ws.Range("E2").Value = 1084.062
$ws.Range("F2").Value = 34.395
$ws.Range("E3").Value = 1002.025
$ws.Range("F3").Value = 35.127
$ws.Range("E4").Value = 1081.163
$ws.Range("F4").Value = 34.709
$ws.Range("E5").Value = 1053.321
$ws.Range("F5").Value = 34.573
$ws.Range("E6").Value = 986.223
$ws.Range("F6").Value = 34.907
$ws.Range("E7").Value = 1068.248
$ws.Range("F7").Value = 35.063
$ws.Range("E8").Value = 1143.015
$ws.Range("F8").Value = 34.641
$ws.Range("E9").Value = 1229.277
$ws.Range("F9").Value = 34.722
$ws.Range("E10").Value = 1087.439
$ws.Range("F10").Value = 30.123
$ws.Range("E11").Value = 1084.954
$ws.Range("F11").Value = 30.877
$ws.Range("E12").Value = 1123.228
$ws.Range("F12").Value = 32.751
$ws.Range("E13").Value = 1106.721
$ws.Range("F13").Value = 35.083
$ws.Range("E14").Value = 971.568
$ws.Range("F14").Value = 34.295
$ws.Range("E15").Value = 1191.449
$ws.Range("F15").Value = 34.896
$ws.Range("E16").Value = 965.316
$ws.Range("F16").Value = 34.655
$ws.Range("E17").Value = 1054.102
$ws.Range("F17").Value = 30.669
$ws.Range("E18").Value = 1167.278
$ws.Range("F18").Value = 29.363
$ws.Range("E19").Value = 933.706
$ws.Range("F19").Value = 30.473
$ws.Range("E20").Value = 1066.96
$ws.Range("F20").Value = 30.997
$ws.Range("E21").Value = 1116.473
$ws.Range("F21").Value = 31.897
$ws.Range("E22").Value = 1171.724
$ws.Range("F22").Value = 30.879
$ws.Range("E23").Value = 1424.404
$ws.Range("F23").Value = 34.422
$ws.Range("E24").Value = 916.125
$ws.Range("F24").Value = 31.009
$ws.Range("E25").Value = 1082.054
$ws.Range("F25").Value = 34.603
$ws.Range("E26").Value = 1026.389
$ws.Range("F26").Value = 31.054
$ws.Range("E27").Value = 1023.598
$ws.Range("F27").Value = 26.652
$ws.Range("E28").Value = 1015.448
$ws.Range("F28").Value = 30.234
$ws.Range("E29").Value = 1156.026
$ws.Range("F29").Value = 29.792
$ws.Range("E30").Value = 994.597
$ws.Range("F30").Value = 30.705
$ws.Range("E31").Value = 1249.841
$ws.Range("F31").Value = 32.179
$ws.Range("E32").Value = 1381.203
$ws.Range("F32").Value = 29.941
$ws.Range("E33").Value = 1398.54
$ws.Range("F33").Value = 32.461
$ws.Range("E34").Value = 910.004
$ws.Range("F34").Value = 30.385
$ws.Range("E35").Value = 977.5309999999999
$ws.Range("F35").Value = 34.06
$ws.Range("E36").Value = 1117.798
$ws.Range("F36").Value = 33.402
$ws.Range("E37").Value = 1046.693
$ws.Range("F37").Value = 33.347
$ws.Range("E38").Value = 1041.886
$ws.Range("F38").Value = 34.191
$ws.Range("E39").Value = 1090.118
$ws.Range("F39").Value = 33.842
$ws.Range("E40").Value = 1103.062
$ws.Range("F40").Value = 34.297
$ws.Range("E41").Value = 988.217
$ws.Range("F41").Value = 34.135
$ws.Range("E42").Value = 979.915
$ws.Range("F42").Value = 34.334
$ws.Range("E43").Value = 1020.68
$ws.Range("F43").Value = 31.375
$ws.Range("E44").Value = 1018.988
$ws.Range("F44").Value = 31.759
$ws.Range("E45").Value = 1332.386
$ws.Range("F45").Value = 27.28
$ws.Range("E46").Value = 1012.249
$ws.Range("F46").Value = 30.572
$ws.Range("E47").Value = 941.836
$ws.Range("F47").Value = 30.668
$ws.Range("E48").Value = 1026.975
$ws.Range("F48").Value = 33.756
$ws.Range("E49").Value = 1170.773
$ws.Range("F49").Value = 29.284
$ws.Range("E50").Value = 1322.373
$ws.Range("F50").Value = 30.118
$ws.Range("E51").Value = 1305.003
$ws.Range("F51").Value = 31.761
$ws.Range("E52").Value = 1044.689
$ws.Range("F52").Value = 30.545
$ws.Range("E53").Value = 1030.866
$ws.Range("F53").Value = 33.282
$ws.Range("E54").Value = 1256.453
$ws.Range("F54").Value = 33.187
$ws.Range("E55").Value = 1004.054
$ws.Range("F55").Value = 30.713
$ws.Range("E56").Value = 969.98
$ws.Range("F56").Value = 34.442
$ws.Range("E57").Value = 975.877
$ws.Range("F57").Value = 33.597
$ws.Range("E58").Value = 936.898
$ws.Range("F58").Value = 33.577

# New Total row 59
$ws.Range("A59").Value = "Total"
$ws.Range("B59").Value = 14042
$ws.Range("C59").Value = 0.5151021008403361
$ws.Range("D59").Value = 160
$ws.Range("E59").Value = 1113.612968380573
$ws.Range("F59").Value = 31.50770801880074
